$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as plain text, preventing Excel from auto-converting
# number-like strings (e.g. "249.41", "37.122.58") into numeric/date values,
# while keeping the cell style/number format unchanged (General, no quote-prefix marker).
function Set-TextValue($range, [string]$value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '37.122.58'
$ws.Range('E2').Value = '  +0.39%  '
Set-TextValue $ws.Range('D3') '2.051.74'
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue $ws.Range('D5') '249.41'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('E6').Value = '  -0.19%  '
Set-TextValue $ws.Range('D7') '60.11'
$ws.Range('E7').Value = '  +8.60%  '
$ws.Range('E8').Value = '  +0.00%  '
Set-TextValue $ws.Range('D9') '0.388'
$ws.Range('E9').Value = '  +1.92%  '
$ws.Range('E10').Value = '  -1.52%  '
$ws.Range('E11').Value = '  +1.79%  '
$ws.Range('E12').Value = '  +7.58%  '
Set-TextValue $ws.Range('D13') '2.349.98'
$ws.Range('E13').Value = '  -0.38%  '
Set-TextValue $ws.Range('D14') '0.835'
$ws.Range('E14').Value = '  +2.53%  '
Set-TextValue $ws.Range('D15') '5.80'
$ws.Range('E15').Value = '  +9.83%  '
Set-TextValue $ws.Range('D16') '2.051.00'
$ws.Range('E16').Value = '  -0.34%  '
Set-TextValue $ws.Range('D17') '18.25'
$ws.Range('E17').Value = '  +28.01%  '
Set-TextValue $ws.Range('D18') '37.146.95'
$ws.Range('E18').Value = '  +0.65%  '
Set-TextValue $ws.Range('D19') '76.20'
$ws.Range('E19').Value = '  +3.46%  '
$ws.Range('E20').Value = '  -3.35%  '
$ws.Range('E21').Value = '  +1.14%  '
Set-TextValue $ws.Range('D22') '238.56'
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('E24').Value = '  -0.29%  '
Set-TextValue $ws.Range('D25') '2.22'
$ws.Range('E25').Value = '  +11.49%  '
$ws.Range('E26').Value = '  +3.59%  '
Set-TextValue $ws.Range('D27') '169.13'
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('E28').Value = '  +0.83%  '
$ws.Range('E29').Value = '  +0.87%  '
$ws.Range('E30').Value = '  +7.94%  '
Set-TextValue $ws.Range('D31') '4.83'
$ws.Range('E31').Value = '  +5.41%  '
Set-TextValue $ws.Range('D32') '0.0630'
$ws.Range('E32').Value = '  +0.40%  '
Set-TextValue $ws.Range('D33') '4.63'
$ws.Range('E33').Value = '  +5.82%  '
Set-TextValue $ws.Range('D34') '0.0891'
$ws.Range('E34').Value = '  +1.42%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('E36').Value = '  -1.13%  '
$ws.Range('E37').Value = '  -1.45%  '
Set-TextValue $ws.Range('D38') '0.109'
$ws.Range('E38').Value = '  +2.77%  '
$ws.Range('E39').Value = '  -0.23%  '
Set-TextValue $ws.Range('D40') '3.18'
$ws.Range('E40').Value = '  +13.85%  '
Set-TextValue $ws.Range('D41') '5.18'
$ws.Range('E41').Value = '  +19.43%  '
$ws.Range('E42').Value = '  +0.09%  '
Set-TextValue $ws.Range('D43') '17.68'
$ws.Range('E43').Value = '  -1.43%  '
$ws.Range('E44').Value = '  -0.16%  '
Set-TextValue $ws.Range('D45') '97.37'
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('E46').Value = '  +3.37%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D47') '1.294.92'
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('B48').Value = 'FTXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range('D48') '3.83'
$ws.Range('E48').Value = '  -7.08%  '
$ws.Range('E49').Value = '  -1.18%  '
$ws.Range('E50').Value = '  -0.55%  '
Set-TextValue $ws.Range('D51') '2.241.71'
$ws.Range('E51').Value = '  -0.21%  '
